# Updates crypto price/volume figures on Sheet1 (columns D "Price" and
# E "Volume(1h)") to the latest scrape values.
#
# The source values are plain text (e.g. "254.16", "3.28%") rather than
# numeric cells, so each target cell is first marked as Text ("@" number
# format) before the new string is written. Without this, Excel's COM
# layer auto-converts numeric-looking strings (and "12.34%"-style strings)
# into real Number cells, which would silently change the cell type and
# the underlying stored value (e.g. "3.28%" -> 0.0328).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-TextValue 'D2' '254.16'
Set-TextValue 'E2' '3.28%'
Set-TextValue 'D3' '27.94'
Set-TextValue 'E3' '-5.48%'
Set-TextValue 'D4' '5.326'
Set-TextValue 'E4' '3.36%'
Set-TextValue 'D5' '0.05843'
Set-TextValue 'E5' '0.85%'
Set-TextValue 'D6' '6.705'
Set-TextValue 'E6' '0.79%'
Set-TextValue 'D7' '0.8658'
Set-TextValue 'E7' '1.63%'
Set-TextValue 'D8' '0.9158'
Set-TextValue 'E8' '5.84%'
Set-TextValue 'D9' '0.1429'
Set-TextValue 'E9' '3.80%'
Set-TextValue 'D10' '0.07178'
Set-TextValue 'E10' '1.32%'
Set-TextValue 'E11' '-1.66%'
Set-TextValue 'D12' '0.09231'
Set-TextValue 'E12' '-1.54%'
Set-TextValue 'D13' '0.001540'
Set-TextValue 'D14' '0.0006047'
Set-TextValue 'E14' '0.58%'
Set-TextValue 'D15' '0.006014'
Set-TextValue 'E15' '-2.02%'
Set-TextValue 'D16' '3.497'
Set-TextValue 'E17' '1.11%'
Set-TextValue 'D18' '2.226'
Set-TextValue 'E18' '0.31%'
Set-TextValue 'E19' '-0.96%'
Set-TextValue 'D20' '0.03455'
Set-TextValue 'E20' '2.67%'
Set-TextValue 'D21' '0.1309'
Set-TextValue 'E21' '2.16%'
Set-TextValue 'D22' '3.528'
Set-TextValue 'E22' '6.20%'
Set-TextValue 'D23' '0.04153'
Set-TextValue 'E23' '0.31%'
Set-TextValue 'D24' '0.1378'
Set-TextValue 'E24' '0.00%'
Set-TextValue 'D25' '0.005110'
Set-TextValue 'D26' '0.001223'
Set-TextValue 'E26' '-0.21%'
Set-TextValue 'D27' '0.0001200'
Set-TextValue 'E27' '-0.78%'
Set-TextValue 'E28' '34.21%'
Set-TextValue 'D40' '0.03859'
Set-TextValue 'E40' '3.00%'
Set-TextValue 'D41' '0.1099'
Set-TextValue 'E41' '2.65%'
Set-TextValue 'D42' '0.002390'
Set-TextValue 'E42' '8.68%'
Set-TextValue 'D43' '0.002950'
Set-TextValue 'E43' '-49.29%'
Set-TextValue 'D44' '0.01095'
Set-TextValue 'E44' '19.23%'
Set-TextValue 'D45' '0.00005244'
Set-TextValue 'E45' '-0.62%'
Set-TextValue 'E46' '0.02%'
Set-TextValue 'D47' '0.08981'
Set-TextValue 'E47' '55.05%'
Set-TextValue 'D48' '0.002154'
Set-TextValue 'E48' '-0.92%'
Set-TextValue 'E49' '0.02%'
Set-TextValue 'E50' '0.02%'
